$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) text labels for signup page fields
# (order matches the new shared-string append order: Birth* then Passwd*)
$ws.Range("G1").Value = "BirthMonth"
$ws.Range("H1").Value = "BirthDay"
$ws.Range("I1").Value = "BirthYear"
$ws.Range("E1").Value = "Passwd"
$ws.Range("F1").Value = "PasswdConfirm"
$ws.Range("O1").Value = "BirthMonthMesEr"
$ws.Range("P1").Value = "BirthDayMesEr"
$ws.Range("Q1").Value = "BirthYearMesEr"
$ws.Range("M1").Value = "PasswdMesEr"
$ws.Range("N1").Value = "PasswdConfirmMesEr"

# Update the active cell selection to A15
$ws.Range("A15").Select()
